# Update Backlog for Task
# Sprint 1 sheet: mark "Design Admin dashboard page" as Finished (with Day2/Day3
# actuals filled in) and move "Design site buider page" to In Process (with a
# Day3 actual logged). Finally leave the selection on K12, matching the saved
# workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1")

# Row 8 - "Design Admin dashboard page"
$ws.Range("F8").Value = "Finished"
$ws.Range("I8").Value = 2
$ws.Range("J8").Value = 4

# Row 9 - "Design site buider page"
$ws.Range("F9").Value = "In Process"
$ws.Range("J9").Value = 1

# Leave the cursor on K12, as captured in the saved file.
$ws.Range("K12").Select() | Out-Null
